$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New product rows appended to the catalog (Enzim SoftShell Mont + Slim Fit kot Pantolon) ---
# Columns: A=urun_adi, B=fiyat, C=kategori, D=gorsel, E=aciklama, F=stok

$softShellDesc = "Su ve rüzgar geçirmez özelliği ile her türlü hava koşulunda konforlu bir kullanım sunar.Polar ve astar detayları ile ekstra sıcaklık ve konfor sağlar.S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$slimDornierDoveDesc = "%98 pamuk ve %2 spandex içeriği sayesinde nefes alabilirlik ve esneklik sağlar, gün boyu konforlu bir kullanım sunar.Slim silueti vücut hatlarınıza mükemmel uyum sağlayarak şık bir profil çizer.31-32-33-34-36-38 Beden seçeneği mevcuttur.Ürünümüz tekli olarak satın alınabilir.Belirtilen fiyatlar adet fiyatıdır."

# Row 53: Enzim SoftShell Mont Füme
$ws.Range("A53").Value = "Enzim SoftShell Mont Füme"
$ws.Range("B53").Value = "600 Tl"
$ws.Range("C53").Value = "Mont"
$ws.Range("D53").Value = "SOFTFÜME.jpg"
$ws.Range("E53").Value = $softShellDesc

# Row 54: Enzim SoftShell Mont Haki
$ws.Range("A54").Value = "Enzim SoftShell Mont Haki"
$ws.Range("B54").Value = "600 Tl"
$ws.Range("C54").Value = "Mont"
$ws.Range("D54").Value = "SOFTHAKİ.jpg"
$ws.Range("E54").Value = $softShellDesc

# Row 55: Enzim SoftShell Mont Siyah
$ws.Range("A55").Value = "Enzim SoftShell Mont Siyah"
$ws.Range("B55").Value = "600 Tl"
$ws.Range("C55").Value = "Mont"
$ws.Range("D55").Value = "SOFTSİYAH.jpg"
$ws.Range("E55").Value = $softShellDesc

# Row 56: Slim Fit kot Pantolon Dornier
$ws.Range("A56").Value = "Slim Fit kot Pantolon Dornier"
$ws.Range("B56").Value = "320 Tl"
$ws.Range("C56").Value = "Jeans"
$ws.Range("D56").Value = "DORNİER.jpg"

# Row 57: Slim Fit kot Pantolon Dove
$ws.Range("A57").Value = "Slim Fit kot Pantolon Dove"
$ws.Range("B57").Value = "320 Tl"
$ws.Range("C57").Value = "Jeans"
$ws.Range("D57").Value = "DOVE.jpg"

# Description filled in afterwards for the Dornier / Dove rows
$ws.Range("E56").Value = $slimDornierDoveDesc
$ws.Range("E57").Value = $slimDornierDoveDesc

# --- "stok" (F) column filled for the existing Kanvas rows and all the new rows ---
$ws.Range("F49").Value = "Var"
$ws.Range("F50").Value = "Var"
$ws.Range("F51").Value = "Var"
$ws.Range("F52").Value = "Var"
$ws.Range("F53").Value = "Var"
$ws.Range("F54").Value = "Var"
$ws.Range("F55").Value = "Var"
$ws.Range("F56").Value = "Var"
$ws.Range("F57").Value = "Var"
$ws.Range("F58").Value = "Var"

# Leave the final selection on E56, matching the saved workbook's cursor position
$ws.Range("E56").Select()
